$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The player roster rows (A2:C17) were re-sorted into a new order.
# Same 16 players/positions/teams, just shuffled rows. Write the new
# order back into the existing cells.
$rows = @(
    ,@("Dennis Schröder",          "PG",         "Brooklyn Nets")
    ,@("Jordan Poole",             "PG,SG",      "Washington Wizards")
    ,@("Mike Conley",              "PG",         "Minnesota Timberwolves")
    ,@("Zach LaVine",              "SG,SF",      "Chicago Bulls")
    ,@("Lauri Markkanen",          "SF,PF",      "Utah Jazz")
    ,@("John Collins",             "PF,C",       "Utah Jazz")
    ,@("Joel Embiid",              "C",          "Philadelphia 76ers")
    ,@("Keyonte George",           "PG,SG",      "Utah Jazz")
    ,@("Kyrie Irving",             "PG,SG",      "Dallas Mavericks")
    ,@("Brandon Boston Jr.",       "SG,SF,PF",   "New Orleans Pelicans")
    ,@("Jalen Williams",           "SG,SF,PF,C", "Oklahoma City Thunder")
    ,@("Jimmy Butler",             "SF,PF",      "Miami Heat")
    ,@("CJ McCollum",              "PG,SG",      "New Orleans Pelicans")
    ,@("RJ Barrett",               "SF,PF",      "Toronto Raptors")
    ,@("Shai Gilgeous-Alexander",  "PG",         "Oklahoma City Thunder")
    ,@("Tobias Harris",            "SF,PF",      "Detroit Pistons")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
